# Typo fixes to 2.doc
#
# This script applies a series of small wording/typo corrections to the
# document using Find/Replace, mirroring the target revision:
#  1. "False alarm caused incorrect reporting ..." ->
#     "False alarms caused by incorrect reporting ..."
#  2. "System must be able to detect and recover ..." ->
#     "System must be able to detect, report  and recover ..."
#  3. "... operations.  An essential part ..." ->
#     "... operations.  A crucial part ..."
#  4. "were identified already in the early design stages" ->
#     "were already identified in the early design stages"
#  5. "those plans should cover on high level mitigation" ->
#     "those plans should cover high level mitigation"
#  6. "should be clearly communicates." ->
#     "should be clearly communicated."

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "False alarm caused incorrect reporting of the system components",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "False alarms caused by incorrect reporting of the system components",
    2) | Out-Null

$d.Content.Find.Execute(
    "System must be able to detect and recover from these kinds of events.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "System must be able to detect, report  and recover from these kinds of events.",
    2) | Out-Null

$d.Content.Find.Execute(
    "must be clearly defined on every level and for all operations.  An essential part ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "must be clearly defined on every level and for all operations.  A crucial part ",
    2) | Out-Null

$d.Content.Find.Execute(
    "were identified already in the early design stages",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "were already identified in the early design stages",
    2) | Out-Null

$d.Content.Find.Execute(
    "those plans should cover on high level mitigation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "those plans should cover high level mitigation",
    2) | Out-Null

$d.Content.Find.Execute(
    "should be clearly communicates.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "should be clearly communicated.",
    2) | Out-Null
